$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing existing rows (15..129) down to (16..130)
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly price record
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 45061
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100114007
$ws.Range("G15").Value = "Jengibre"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("N15").Value = "$/caja 13 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1346
$ws.Range("Q15").Value = 13
$ws.Range("R15").Value = "Hortaliza"
